$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v_A2 = @"
YEZA INN
Opened in 2025
8.5/10
Very good2 reviews
Near Yaya CentreShow on Map
2-bed Room
x4
Free Cancellation
Breakfast included
Only 1 left at this price
Special Discount
10% off
₹ 6,160
₹ 5,487
Total price: ₹ 32,429
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices
"@
$ws.Range("A2").Value2 = $v_A2

$v_B2 = @"
Woodmere Serviced Apartment
8.6/10
Very good70 reviews
Near Yaya CentreShow on Map
Standard Two-Bedroom Apartment
x4
Entire unit 59㎡
2 bedrooms
2 beds
₹ 5,237
Total price: ₹ 30,869
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("B2").Value2 = $v_B2

$v_C2 = @"
Holiday Inn NAIROBI TWO RIVERS MALL by IHG
9.7/10
Amazing94 reviews
"Great location"
"Great service"
Near Village MarketShow on Map
No. 13 of 4-Star Select Hotels in Nairobi
Family Room
x4
Free Cancellation
Last booked 19 hrs ago
₹ 19,978
Total price: ₹ 127,861
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("C2").Value2 = $v_C2

$v_A3 = @"
Javilla Eagles Safari Guest house
7.7/10
Good34 reviews
8.4 km from centreShow on Map
Comfort Apartment, 2 Bedrooms, Non Smoking, Ground Floor
x4
Entire unit 55㎡
2 bedrooms
2 beds
Only 1 left at this price
₹ 5,056
Total price: ₹ 29,326
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A3").Value2 = $v_A3

$v_B3 = @"
Kenya Comfort Suites
7.9/10
Good11 reviews
Near Yaya CentreShow on Map
Standard Quadruple Room
x4
₹ 4,783
Total price: ₹ 29,380
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("B3").Value2 = $v_B3

$v_C3 = @"
JW Marriott Hotel Nairobi
9.8/10
Outstanding46 reviews
"Great stay!"
"Great service"
Near National Museum of KenyaShow on Map
No. 1 of Luxury Hotels in Nairobi
3 Bedroom Apartment, Bedroom 1: 1 King, Bedroom 2: 1 King, Bedroom 3: 2 Doubles
x4
Free Cancellation
Breakfast included
Last booked 15 hrs ago
₹ 203,323
Total price: ₹ 1,283,575
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("C3").Value2 = $v_C3

$v_A4 = @"
Woodmere Serviced Apartment
8.6/10
Very good70 reviews
Near Yaya CentreShow on Map
Standard Two-Bedroom Apartment
x4
Entire unit 59㎡
2 bedrooms
2 beds
₹ 5,237
Total price: ₹ 30,869
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A4").Value2 = $v_A4

$v_B4 = @"
Eldon Apartments & Suites
7.1/10
5 reviews
Near Wilson AirportShow on Map
Standard Two-Bedroom Apartment
x4
Entire apartment 18㎡
2 bedrooms
3 beds
Special Discount
₹ 111 off
₹ 4,439
₹ 4,328
Total price: ₹ 25,544
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("B4").Value2 = $v_B4

$v_C4 = @"
Yaya Hotel & Apartments
9.4/10
Amazing49 reviews
"Clean and tidy"
"Great service"
Near Yaya CentreShow on Map
No. 19 of 4-Star Select Hotels in Nairobi
Two-Bedroom Luxury Apartment
x4
Entire unit 125㎡
2 bedrooms
2 beds
Free Cancellation
Earn ₹ 961.74 in Trip Coins
Last booked 5 hrs ago
Limited Time Offer
15% off
₹ 19,625
₹ 16,298
Total price: ₹ 96,160
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices
"@
$ws.Range("C4").Value2 = $v_C4

$v_A5 = @"
Eldon Apartments & Suites
7.1/10
5 reviews
Near Wilson AirportShow on Map
Standard Two-Bedroom Apartment
x4
Entire apartment 18㎡
2 bedrooms
3 beds
Special Discount
₹ 111 off
₹ 4,439
₹ 4,328
Total price: ₹ 25,544
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A5").Value2 = $v_A5

$v_B5 = @"
The King Post
7.7/10
Good17 reviews
Near The Sarit Expo CentreShow on Map
Three-Bedroom Apartment
x6
Entire unit
3 bedrooms
4 beds
₹ 4,458
Total price: ₹ 26,254
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("B5").Value2 = $v_B5

$v_C5 = @"
Pan Pacific Serviced Suites Nairobi
9.4/10
Amazing48 reviews
"Great rooms"
"Great location"
Near National Museum of KenyaShow on Map
No. 13 of Premium Hotels in Nairobi
Two Bedroom Suite King & Twin
x4
Entire unit 108㎡
2 bedrooms
2 beds
Last booked 2 hrs ago
₹ 23,226
Total price: ₹ 145,164
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("C5").Value2 = $v_C5

$v_A6 = @"
Kenya Comfort Suites
7.9/10
Good11 reviews
Near Yaya CentreShow on Map
Standard Quadruple Room
x4
₹ 4,783
Total price: ₹ 29,380
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A6").Value2 = $v_A6

$v_B6 = @"
Arcadia Hotel
Renovated in 2025
8.8/10
Very good13 reviews
Near Yaya CentreShow on Map
Two-Bedroom Suite
x4
Entire unit 110㎡
2 bedrooms
2 beds
Free Cancellation
Only 5 left at this price
₹ 7,701
Total price: ₹ 44,666
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices
"@
$ws.Range("B6").Value2 = $v_B6

$v_C6 = @"
Palacina the Residence & the Suites
9.2/10
Great48 reviews
Near Yaya CentreShow on Map
No. 7 of Premium Hotels in Nairobi
2 Bedroom Executive Penthouse
x4
Entire unit 168㎡
2 bedrooms
3 beds
Free Cancellation
Breakfast included
Earn ₹ 2,505.41 in Trip Coins
Only 1 left at this price
Special Discount
20% off
₹ 36,660
₹ 28,309
Total price: ₹ 167,023
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices
"@
$ws.Range("C6").Value2 = $v_C6

$v_A7 = @"
Holiday Inn NAIROBI TWO RIVERS MALL by IHG
Ad
9.7/10
Amazing94 reviews
"Great location"
"Great service"
Near Village MarketShow on Map
No. 13 of 4-Star Select Hotels in Nairobi
Family Room
x4
Free Cancellation
Last booked 19 hrs ago
₹ 19,978
Total price: ₹ 127,861
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A7").Value2 = $v_A7

$v_B7 = @"
Holiday Inn NAIROBI TWO RIVERS MALL by IHG
Ad
9.7/10
Amazing94 reviews
"Great location"
"Great service"
Near Village MarketShow on Map
No. 13 of 4-Star Select Hotels in Nairobi
Family Room
x4
Free Cancellation
Last booked 19 hrs ago
₹ 19,978
Total price: ₹ 127,861
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("B7").Value2 = $v_B7

$v_C7 = @"
Executive Residency by Best Western Nairobi
9.2/10
Great46 reviews
Near The Sarit Expo CentreShow on Map
Two Bedroom Apartment
x4
Breakfast included
₹ 17,043
Total price: ₹ 104,816
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("C7").Value2 = $v_C7

$v_A8 = @"
The King Post
7.7/10
Good17 reviews
Near The Sarit Expo CentreShow on Map
Three-Bedroom Apartment
x6
Entire unit
3 bedrooms
4 beds
₹ 4,458
Total price: ₹ 26,254
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A8").Value2 = $v_A8

$v_B8 = @"
Lux Suites Riara One Residency Angama
Opened in 2025
9.9/10
Outstanding17 reviews
"Clean and tidy"
"Great location"
Near Yaya CentreShow on Map
Family Room
x4
Entire apartment 98㎡
2 bedrooms
2 beds
Free Cancellation
Breakfast included
Earn ₹ 1,322.62 in Trip Coins
Only 5 left at this price
Special Discount
8% off
₹ 16,393
₹ 14,942
Total price: ₹ 88,151
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices
"@
$ws.Range("B8").Value2 = $v_B8

$v_C8 = @"
Mövenpick Hotel & Residences Nairobi
9.0/10
Great97 reviews
"Great service"
"Great location"
Near The Sarit Expo CentreShow on Map
No. 1 of Gourmet Hotels in Nairobi
Two-Bedroom Residence
x4
Free Cancellation
Last booked 22 hrs ago
₹ 18,676
Total price: ₹ 116,723
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("C8").Value2 = $v_C8

$v_A9 = @"
Kester International Apartment Hotel
Opened in 2025
9.8/10
Outstanding28 reviews
"Great service"
"Great rooms"
Near Yaya CentreShow on Map
Boutique 2-bedroom And 1-living Room Suite
x4
Entire apartment 95㎡
2 bedrooms
3 beds
Earn ₹ 546.73 in Trip Coins
Only 3 left at this price
Special Discount
11% off
₹ 7,139
₹ 6,281
Total price: ₹ 36,432
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices
"@
$ws.Range("A9").Value2 = $v_A9

$v_B9 = @"
Lavington Residences By Trianum
Opened in 2025
9.6/10
Amazing7 reviews
Near Yaya CentreShow on Map
Executive Two-Bedroom Apartment
x4
Only 1 left at this price
₹ 9,559
Total price: ₹ 56,396
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices
"@
$ws.Range("B9").Value2 = $v_B9

$v_C9 = @"
Fairview Hotel Nairobi, Vignette Collection by IHG
Renovated in 2025
9.0/10
Great49 reviews
Near Giraffe manorShow on Map
No. 17 of 4-Star Select Hotels in Nairobi
Fairview Suite
x4
Free Cancellation
₹ 54,628
Total price: ₹ 349,616
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("C9").Value2 = $v_C9

$v_A10 = @"
Maskan Suites
8.9/10
Very good34 reviews
Near Yaya CentreShow on Map
Superior Apartment, 2 Bedrooms, Private Bathroom, City View
x4
Entire apartment
2 bedrooms
2 beds
Free Cancellation
₹ 6,943
Total price: ₹ 40,926
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A10").Value2 = $v_A10

$v_B10 = @"
Mövenpick Hotel & Residences Nairobi
9.0/10
Great97 reviews
"Great service"
"Great location"
Near The Sarit Expo CentreShow on Map
No. 1 of Gourmet Hotels in Nairobi
Two-Bedroom Residence
x4
Free Cancellation
Last booked 22 hrs ago
₹ 18,676
Total price: ₹ 116,723
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("B10").Value2 = $v_B10

$v_C10 = @"
Windsor Golf Hotel & Country Club
New to Trip.com
8.8/10
Very good64 reviews
9.7 km from centreShow on Map
No. 6 of 4-Star Select Hotels in Nairobi
Two- Bedroom Cottage
x4
Entire unit 65㎡
1 bedroom
4 beds
Breakfast included
₹ 36,881
Total price: ₹ 217,596
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("C10").Value2 = $v_C10

$v_A11 = @"
Mövenpick Hotel & Residences Nairobi
9.0/10
Great97 reviews
"Great service"
"Great location"
Near The Sarit Expo CentreShow on Map
No. 1 of Gourmet Hotels in Nairobi
Two-Bedroom Residence
x4
Free Cancellation
Last booked 22 hrs ago
₹ 18,676
Total price: ₹ 116,723
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A11").Value2 = $v_A11

$v_B11 = @"
Yaya Hotel & Apartments
9.4/10
Amazing49 reviews
"Clean and tidy"
"Great service"
Near Yaya CentreShow on Map
No. 19 of 4-Star Select Hotels in Nairobi
Two-Bedroom Luxury Apartment
x4
Entire unit 125㎡
2 bedrooms
2 beds
Free Cancellation
Earn ₹ 961.74 in Trip Coins
Last booked 5 hrs ago
Limited Time Offer
15% off
₹ 19,625
₹ 16,298
Total price: ₹ 96,160
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices
"@
$ws.Range("B11").Value2 = $v_B11

$v_C11 = @"
Four Points by Sheraton Nairobi Hurlingham
8.8/10
Very good60 reviews
"Great service"
"Delicious breakfast"
Near Giraffe manorShow on Map
Executive Suite
x4
₹ 29,249
Total price: ₹ 187,235
1 room × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("C11").Value2 = $v_C11

$v_A12 = @"
Mercure Nairobi Upper Hill
Ad
8.9/10
Very good60 reviews
Near Giraffe manorShow on Map
2
Classic Superior King Room
Free Cancellation
₹ 11,333
Total price: ₹ 141,667
2 rooms × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("A12").Value2 = $v_A12

$v_B12 = @"
Mercure Nairobi Upper Hill
Ad
8.9/10
Very good60 reviews
Near Giraffe manorShow on Map
2
Classic Superior King Room
Free Cancellation
₹ 11,333
Total price: ₹ 141,667
2 rooms × 5 nights incl. taxes & fees
Check Availability
"@
$ws.Range("B12").Value2 = $v_B12
